$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.883.46'
$ws.Range("E2").Value = '  +5.51%  '
$ws.Range("D3").Value = '2.417.00'
$ws.Range("E3").Value = '  +2.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +1.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.43'
$ws.Range("E5").Value = '  +2.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.43'
$ws.Range("E6").Value = '  +5.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.540'
$ws.Range("E8").Value = '  +2.32%  '
$ws.Range("D9").Value = '2.447.15'
$ws.Range("E9").Value = '  +3.60%  '
$ws.Range("E10").Value = '  +5.99%  '
$ws.Range("E11").Value = '  +1.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.22'
$ws.Range("E12").Value = '  +3.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.353'
$ws.Range("E13").Value = '  +5.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.04'
$ws.Range("E14").Value = '  +6.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000177'
$ws.Range("E15").Value = '  +7.86%  '
$ws.Range("D16").Value = '2.963.49'
$ws.Range("E16").Value = '  +7.22%  '
$ws.Range("D17").Value = '62.975.74'
$ws.Range("E17").Value = '  +5.94%  '
$ws.Range("D18").Value = '2.454.30'
$ws.Range("E18").Value = '  +4.06%  '
$ws.Range("E19").Value = '  -1.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.94'
$ws.Range("E20").Value = '  +4.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '327.85'
$ws.Range("E21").Value = '  +2.30%  '
$ws.Range("E22").Value = '  +2.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.04'
$ws.Range("E23").Value = '  +12.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.63'
$ws.Range("E25").Value = '  +2.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '620.70'
$ws.Range("E26").Value = '  +12.14%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.39'
$ws.Range("E27").Value = '  +3.60%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0980'
$ws.Range("E28").Value = '  +7.17%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '2.565.37'
$ws.Range("E29").Value = '  +3.82%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.14'
$ws.Range("E30").Value = '  +2.18%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.41'
$ws.Range("E31").Value = '  +8.22%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.85'
$ws.Range("E32").Value = '  +4.41%  '
$ws.Range("E33").Value = '  +5.82%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.48'
$ws.Range("E34").Value = '  +4.88%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.995'
$ws.Range("E35").Value = '  -0.46%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.75'
$ws.Range("E36").Value = '  +5.11%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.373'
$ws.Range("E37").Value = '  +2.05%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.67'
$ws.Range("E38").Value = '  +1.09%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.38'
$ws.Range("E39").Value = '  +7.99%  '
$ws.Range("B40").Value = 'EthereumClassic'
$ws.Range("C40").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.65'
$ws.Range("E40").Value = '  +3.23%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.76'
$ws.Range("E41").Value = '  +15.72%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.76'
$ws.Range("E42").Value = '  +7.26%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.37'
$ws.Range("E43").Value = '  +2.62%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.0₆0285'
$ws.Range("E45").Value = '  -3.28%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '144.54'
$ws.Range("E46").Value = '  +5.00%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.59'
$ws.Range("E47").Value = '  +2.42%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.23'
$ws.Range("E48").Value = '  +6.14%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.600'
$ws.Range("E49").Value = '  +2.87%  '
$ws.Range("B50").Value = 'Hedera'
$ws.Range("C50").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0516'
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0917'
$ws.Range("E51").Value = '  +2.96%  '
